$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 133 ("photograph" / 写真) which duplicates row 129.
# This shifts rows 134:146 up to 133:145.
$ws.Rows.Item(133).Delete()
